$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743717012329/Ggg0CQFNgb-camera-video-segment-driver-1743717014829.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNKXJ7YWB%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T210018Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKH%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDjcp2%2F1XPHqFBU8Ha6e0%2BboFHzqMMuSizFvSsvZmG8MwIgdvPD0g1GYwyRrIUpf8wrJ9AAWJWlb1lWyL%2ByOQlXR%2BUq3QMIGhAEGgw3ODEyMDQ5NDIyNDQiDIP%2B%2FDt3WrNLrXs8NCq6AwcOSuMqbzzYfqjLT%2FqfZe0Yk9BySnPW%2FtDgGaswNt2eyLVX%2Bpz7BhS98DBZoTKxHW49%2B63ZiEWys0HxTcEgqvvqePPj7fZL0ljICWfQUlU4Fj2Gegvd3gFcd9XIx8Pe%2B%2BDoVpE%2BggWnNYrLWL4RQ%2FXhO3U%2FIDczwW17eJiwYo1UC3gTy1DVbVUebB4j5YYJ8Bs4zj%2F%2B%2Bcq%2BBVxo7%2BLjEGpW%2B7c2ERKOEYx4DwWks9sg8dQZH1FFXc9URzAxbhoIzWfWSFpK3tYCtvwduzAceq%2Bcrqlf%2BhRl6NPBb4Daru%2FPT2jfQZFeP2geWn7g0qeBXfe4aXMEIMRW9H%2BdyH6sxfvfcItWyX7Wlw8EohO2HX0PSgwAt0nUHhRe2zOC6M5t3c9Rp3mSOYchTkcLV9au0WeotoAZo5dNraL4eKNtJVnWRO2jR7GOWY0v7Nb29zH%2FMBsmIRc9P%2B8fC9QLb4brfyY6BGbwoqyBu6L5pt3toNcffXowdz6uoRkQfow4K7uFENF2z1%2F0%2F65TRpIlxifVwudSqTqKOnmWvmJbJygRM4XP3ggl4CLH%2FJQAo%2FRS8CMMdmplxCxrdnRjPZAw8JvAvwY6pQFF51Rkp2UpcWbbMVJ6DMqL%2FFkVO%2BJkkvfld%2Bj33Z%2B5Xltl%2BpA%2Bu0JJG4R16IDKb4aHm%2BfUUkZbZ%2Fa9yUgsn8Vw08qNMay2NxG99vat2314FXjbE6TfifRlsOttvT0r5OiYNrAMXn5hzpjbYFb2njGQUr%2B5SVXYh9qMs8ZFFvGczXZaR0OPy5kfQgdmYMdhofMNxjVZnFUzv7HxjqkLOXs000%2F8%2BpM%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2005%3A00%3A18%20GMT&X-Amz-Signature=2a9ea820fd2a9784a14765f2fee594d5e38725c23b327164a1db30528c5fcf49'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743713574485/FgwSeOefdQ-camera-video-segment-driver-1743713576985.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNKXJ7YWB%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T210018Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKH%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDjcp2%2F1XPHqFBU8Ha6e0%2BboFHzqMMuSizFvSsvZmG8MwIgdvPD0g1GYwyRrIUpf8wrJ9AAWJWlb1lWyL%2ByOQlXR%2BUq3QMIGhAEGgw3ODEyMDQ5NDIyNDQiDIP%2B%2FDt3WrNLrXs8NCq6AwcOSuMqbzzYfqjLT%2FqfZe0Yk9BySnPW%2FtDgGaswNt2eyLVX%2Bpz7BhS98DBZoTKxHW49%2B63ZiEWys0HxTcEgqvvqePPj7fZL0ljICWfQUlU4Fj2Gegvd3gFcd9XIx8Pe%2B%2BDoVpE%2BggWnNYrLWL4RQ%2FXhO3U%2FIDczwW17eJiwYo1UC3gTy1DVbVUebB4j5YYJ8Bs4zj%2F%2B%2Bcq%2BBVxo7%2BLjEGpW%2B7c2ERKOEYx4DwWks9sg8dQZH1FFXc9URzAxbhoIzWfWSFpK3tYCtvwduzAceq%2Bcrqlf%2BhRl6NPBb4Daru%2FPT2jfQZFeP2geWn7g0qeBXfe4aXMEIMRW9H%2BdyH6sxfvfcItWyX7Wlw8EohO2HX0PSgwAt0nUHhRe2zOC6M5t3c9Rp3mSOYchTkcLV9au0WeotoAZo5dNraL4eKNtJVnWRO2jR7GOWY0v7Nb29zH%2FMBsmIRc9P%2B8fC9QLb4brfyY6BGbwoqyBu6L5pt3toNcffXowdz6uoRkQfow4K7uFENF2z1%2F0%2F65TRpIlxifVwudSqTqKOnmWvmJbJygRM4XP3ggl4CLH%2FJQAo%2FRS8CMMdmplxCxrdnRjPZAw8JvAvwY6pQFF51Rkp2UpcWbbMVJ6DMqL%2FFkVO%2BJkkvfld%2Bj33Z%2B5Xltl%2BpA%2Bu0JJG4R16IDKb4aHm%2BfUUkZbZ%2Fa9yUgsn8Vw08qNMay2NxG99vat2314FXjbE6TfifRlsOttvT0r5OiYNrAMXn5hzpjbYFb2njGQUr%2B5SVXYh9qMs8ZFFvGczXZaR0OPy5kfQgdmYMdhofMNxjVZnFUzv7HxjqkLOXs000%2F8%2BpM%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2005%3A00%3A18%20GMT&X-Amz-Signature=1f394c4c4fdef827ac96b9b335c7a490adc5de1b23e8570e31f565a5cccfffd8'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743712999664/UzzO0YQTJX-camera-video-segment-driver-1743713002164.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNKXJ7YWB%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T210018Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKH%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDjcp2%2F1XPHqFBU8Ha6e0%2BboFHzqMMuSizFvSsvZmG8MwIgdvPD0g1GYwyRrIUpf8wrJ9AAWJWlb1lWyL%2ByOQlXR%2BUq3QMIGhAEGgw3ODEyMDQ5NDIyNDQiDIP%2B%2FDt3WrNLrXs8NCq6AwcOSuMqbzzYfqjLT%2FqfZe0Yk9BySnPW%2FtDgGaswNt2eyLVX%2Bpz7BhS98DBZoTKxHW49%2B63ZiEWys0HxTcEgqvvqePPj7fZL0ljICWfQUlU4Fj2Gegvd3gFcd9XIx8Pe%2B%2BDoVpE%2BggWnNYrLWL4RQ%2FXhO3U%2FIDczwW17eJiwYo1UC3gTy1DVbVUebB4j5YYJ8Bs4zj%2F%2B%2Bcq%2BBVxo7%2BLjEGpW%2B7c2ERKOEYx4DwWks9sg8dQZH1FFXc9URzAxbhoIzWfWSFpK3tYCtvwduzAceq%2Bcrqlf%2BhRl6NPBb4Daru%2FPT2jfQZFeP2geWn7g0qeBXfe4aXMEIMRW9H%2BdyH6sxfvfcItWyX7Wlw8EohO2HX0PSgwAt0nUHhRe2zOC6M5t3c9Rp3mSOYchTkcLV9au0WeotoAZo5dNraL4eKNtJVnWRO2jR7GOWY0v7Nb29zH%2FMBsmIRc9P%2B8fC9QLb4brfyY6BGbwoqyBu6L5pt3toNcffXowdz6uoRkQfow4K7uFENF2z1%2F0%2F65TRpIlxifVwudSqTqKOnmWvmJbJygRM4XP3ggl4CLH%2FJQAo%2FRS8CMMdmplxCxrdnRjPZAw8JvAvwY6pQFF51Rkp2UpcWbbMVJ6DMqL%2FFkVO%2BJkkvfld%2Bj33Z%2B5Xltl%2BpA%2Bu0JJG4R16IDKb4aHm%2BfUUkZbZ%2Fa9yUgsn8Vw08qNMay2NxG99vat2314FXjbE6TfifRlsOttvT0r5OiYNrAMXn5hzpjbYFb2njGQUr%2B5SVXYh9qMs8ZFFvGczXZaR0OPy5kfQgdmYMdhofMNxjVZnFUzv7HxjqkLOXs000%2F8%2BpM%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2005%3A00%3A18%20GMT&X-Amz-Signature=650714211889b64f734df2ac870de53fe69964256ae53a5bec4d8f7818d99748'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743712546260/Jh2p7rn04n-camera-video-segment-driver-1743712548760.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNKXJ7YWB%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T210018Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKH%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDjcp2%2F1XPHqFBU8Ha6e0%2BboFHzqMMuSizFvSsvZmG8MwIgdvPD0g1GYwyRrIUpf8wrJ9AAWJWlb1lWyL%2ByOQlXR%2BUq3QMIGhAEGgw3ODEyMDQ5NDIyNDQiDIP%2B%2FDt3WrNLrXs8NCq6AwcOSuMqbzzYfqjLT%2FqfZe0Yk9BySnPW%2FtDgGaswNt2eyLVX%2Bpz7BhS98DBZoTKxHW49%2B63ZiEWys0HxTcEgqvvqePPj7fZL0ljICWfQUlU4Fj2Gegvd3gFcd9XIx8Pe%2B%2BDoVpE%2BggWnNYrLWL4RQ%2FXhO3U%2FIDczwW17eJiwYo1UC3gTy1DVbVUebB4j5YYJ8Bs4zj%2F%2B%2Bcq%2BBVxo7%2BLjEGpW%2B7c2ERKOEYx4DwWks9sg8dQZH1FFXc9URzAxbhoIzWfWSFpK3tYCtvwduzAceq%2Bcrqlf%2BhRl6NPBb4Daru%2FPT2jfQZFeP2geWn7g0qeBXfe4aXMEIMRW9H%2BdyH6sxfvfcItWyX7Wlw8EohO2HX0PSgwAt0nUHhRe2zOC6M5t3c9Rp3mSOYchTkcLV9au0WeotoAZo5dNraL4eKNtJVnWRO2jR7GOWY0v7Nb29zH%2FMBsmIRc9P%2B8fC9QLb4brfyY6BGbwoqyBu6L5pt3toNcffXowdz6uoRkQfow4K7uFENF2z1%2F0%2F65TRpIlxifVwudSqTqKOnmWvmJbJygRM4XP3ggl4CLH%2FJQAo%2FRS8CMMdmplxCxrdnRjPZAw8JvAvwY6pQFF51Rkp2UpcWbbMVJ6DMqL%2FFkVO%2BJkkvfld%2Bj33Z%2B5Xltl%2BpA%2Bu0JJG4R16IDKb4aHm%2BfUUkZbZ%2Fa9yUgsn8Vw08qNMay2NxG99vat2314FXjbE6TfifRlsOttvT0r5OiYNrAMXn5hzpjbYFb2njGQUr%2B5SVXYh9qMs8ZFFvGczXZaR0OPy5kfQgdmYMdhofMNxjVZnFUzv7HxjqkLOXs000%2F8%2BpM%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2005%3A00%3A18%20GMT&X-Amz-Signature=fa0a6d2eb0d512d5cc612a4d0ad7d3c9379d45cbd52ee67113a42b863db85dff'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743703000742/spTbAbfADq-camera-video-segment-driver-1743703003242.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNKXJ7YWB%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T210018Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKH%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDjcp2%2F1XPHqFBU8Ha6e0%2BboFHzqMMuSizFvSsvZmG8MwIgdvPD0g1GYwyRrIUpf8wrJ9AAWJWlb1lWyL%2ByOQlXR%2BUq3QMIGhAEGgw3ODEyMDQ5NDIyNDQiDIP%2B%2FDt3WrNLrXs8NCq6AwcOSuMqbzzYfqjLT%2FqfZe0Yk9BySnPW%2FtDgGaswNt2eyLVX%2Bpz7BhS98DBZoTKxHW49%2B63ZiEWys0HxTcEgqvvqePPj7fZL0ljICWfQUlU4Fj2Gegvd3gFcd9XIx8Pe%2B%2BDoVpE%2BggWnNYrLWL4RQ%2FXhO3U%2FIDczwW17eJiwYo1UC3gTy1DVbVUebB4j5YYJ8Bs4zj%2F%2B%2Bcq%2BBVxo7%2BLjEGpW%2B7c2ERKOEYx4DwWks9sg8dQZH1FFXc9URzAxbhoIzWfWSFpK3tYCtvwduzAceq%2Bcrqlf%2BhRl6NPBb4Daru%2FPT2jfQZFeP2geWn7g0qeBXfe4aXMEIMRW9H%2BdyH6sxfvfcItWyX7Wlw8EohO2HX0PSgwAt0nUHhRe2zOC6M5t3c9Rp3mSOYchTkcLV9au0WeotoAZo5dNraL4eKNtJVnWRO2jR7GOWY0v7Nb29zH%2FMBsmIRc9P%2B8fC9QLb4brfyY6BGbwoqyBu6L5pt3toNcffXowdz6uoRkQfow4K7uFENF2z1%2F0%2F65TRpIlxifVwudSqTqKOnmWvmJbJygRM4XP3ggl4CLH%2FJQAo%2FRS8CMMdmplxCxrdnRjPZAw8JvAvwY6pQFF51Rkp2UpcWbbMVJ6DMqL%2FFkVO%2BJkkvfld%2Bj33Z%2B5Xltl%2BpA%2Bu0JJG4R16IDKb4aHm%2BfUUkZbZ%2Fa9yUgsn8Vw08qNMay2NxG99vat2314FXjbE6TfifRlsOttvT0r5OiYNrAMXn5hzpjbYFb2njGQUr%2B5SVXYh9qMs8ZFFvGczXZaR0OPy5kfQgdmYMdhofMNxjVZnFUzv7HxjqkLOXs000%2F8%2BpM%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2005%3A00%3A18%20GMT&X-Amz-Signature=f7fec460342fca2401fb304f7dfc24e3539a8ddc4c7641d6f82081052db7ff6a'
